$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix week value for rows 55-57 (7 -> 6)
$ws.Range("C55").Value = 6
$ws.Range("C56").Value = 6
$ws.Range("C57").Value = 6

# Add new rows 58-60 for the bva poll (10/15)
# Row 58
$ws.Range("A58").Value = 17
$ws.Range("B58").Value = 2021
$ws.Range("C58").Value = 6
$ws.Range("D58").Value = 10
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = "bva"
$ws.Range("G58").Value = "online"
$ws.Range("H58").Value = "excluded"
$ws.Range("I58").Value = 876
$ws.Range("J58").Value = 1
$ws.Range("K58").Value = 1
$ws.Range("L58").Value = 8
$ws.Range("M58").Value = 1.5
$ws.Range("N58").Value = 4
$ws.Range("O58").Value = 8
$ws.Range("P58").Value = 4
$ws.Range("Q58").Value = 27
$ws.Range("R58").Value = 10
$ws.Range("U58").Value = 1.5
$ws.Range("V58").Value = 2.5
$ws.Range("W58").Value = 17
$ws.Range("X58").Value = 14
$ws.Range("Y58").Value = 0.5

# Row 59
$ws.Range("A59").Value = 17
$ws.Range("B59").Value = 2021
$ws.Range("C59").Value = 6
$ws.Range("D59").Value = 10
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = "bva"
$ws.Range("G59").Value = "online"
$ws.Range("H59").Value = "excluded"
$ws.Range("I59").Value = 879
$ws.Range("J59").Value = 1
$ws.Range("K59").Value = 1
$ws.Range("L59").Value = 8.5
$ws.Range("M59").Value = 1.5
$ws.Range("N59").Value = 4
$ws.Range("O59").Value = 8
$ws.Range("P59").Value = 4
$ws.Range("Q59").Value = 28
$ws.Range("S59").Value = 8.5
$ws.Range("U59").Value = 1.5
$ws.Range("V59").Value = 3
$ws.Range("W59").Value = 17
$ws.Range("X59").Value = 13
$ws.Range("Y59").Value = 1

# Row 60
$ws.Range("A60").Value = 17
$ws.Range("B60").Value = 2021
$ws.Range("C60").Value = 6
$ws.Range("D60").Value = 10
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = "bva"
$ws.Range("G60").Value = "online"
$ws.Range("H60").Value = "excluded"
$ws.Range("I60").Value = 885
$ws.Range("J60").Value = 0.5
$ws.Range("K60").Value = 1.5
$ws.Range("L60").Value = 8
$ws.Range("M60").Value = 1.5
$ws.Range("N60").Value = 4
$ws.Range("O60").Value = 8
$ws.Range("P60").Value = 4
$ws.Range("Q60").Value = 26
$ws.Range("T60").Value = 12
$ws.Range("U60").Value = 1
$ws.Range("V60").Value = 2.5
$ws.Range("W60").Value = 16
$ws.Range("X60").Value = 14
$ws.Range("Y60").Value = 1

# Update selection to match final cursor location
$ws.Range("I61").Select()
